$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H1 (3 -> 4)
$ws.Range("H1").Value = 4

# Populate new column I as a duplicate of column H (post-update)
$ws.Range("I1").Value = 4
$ws.Range("I2").Value = 2
$ws.Range("I3").Value = 2
$ws.Range("I4").Value = 2
$ws.Range("I5").Value = 2
$ws.Range("I6").Value = 2
$ws.Range("I7").Value = 1
$ws.Range("I8").Value = 2
$ws.Range("I9").Value = 2
$ws.Range("I10").Value = 2
$ws.Range("I11").Value = 2
$ws.Range("I12").Value = 2
$ws.Range("I13").Value = 2
$ws.Range("I14").Value = 2
$ws.Range("I15").Value = 2
$ws.Range("I16").Value = 2
$ws.Range("I17").Value = 2
$ws.Range("I18").Value = 2
$ws.Range("I19").Value = 2
$ws.Range("I20").Value = 2
$ws.Range("I21").Value = 2
$ws.Range("I22").Value = 2
$ws.Range("I23").Value = 2
$ws.Range("I24").Value = 2
$ws.Range("I25").Value = 2
$ws.Range("I26").Value = 2
$ws.Range("I27").Value = 2
